$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SimilarityStatus")

# Update the "algo" column (C) text for methods 6, 8, 9, 10 (rows 20, 22, 23, 24)
$ws.Range("C20").Value = "If the branch has more than 3 children (aka >2 leaves + 1 ambi), check if seq stay with any leave. Else, insert seq to ambi (create 1 if hasn't). Then, if the priority of ambi >= stay threshold, make this ambi node non-ambi."
$ws.Range("C22").Value = "do like method 5 or 7 then see if the new t_parent stay with any of the super (it should be) then if the child is closer to t_parent, create super for once then move the child there. Update all above"
$ws.Range("C23").Value = "do like method 6 or 7 then see if the new t_parent stay with any of the super (it should be) then if the child is closer to t_parent, create super for once then move the child there. Update all above"
$ws.Range("C24").Value = "do like method 7 then see if the new t_parent stay with any of the super (it should be) then if the child is closer to t_parent, create super for once then move the child there. Update all above"

# Auto-fit the row height for row 20 since its text grew substantially
$ws.Rows.Item(20).EntireRow.AutoFit()
$ws.Rows.Item(20).RowHeight = 145

# Update the sheet view: scroll position and selection to reflect where the edit was made
$ws.Range("C20").Select()
$excel.ActiveWindow.ScrollRow = 18

$ws.Activate()
